$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that gets bumped by one day
# for every data row (rows 2 through 307) as part of an automatic update.
for ($row = 2; $row -le 307; $row++) {
    $ws.Cells.Item($row, 3).Value = 45206
}
